{"js": "// Update the date heading and every two-digit-by-two-digit multiplication\n// prompt in the practice-sheet table. Each \"old\" value below occurs exactly\n// once in the document, so a scoped `body.search()` + replace is safe and\n// unambiguous for every pair.\nconst replacements = [\n  [\"2024-02-26 Monday\", \"2024-02-27 Tuesday\"],\n  [\"37\u00d788=\", \"66\u00d788=\"],\n  [\"49\u00d723=\", \"40\u00d721=\"],\n  [\"20\u00d739=\", \"99\u00d724=\"],\n  [\"81\u00d774=\", \"37\u00d737=\"],\n  [\"37\u00d773=\", \"12\u00d762=\"],\n  [\"33\u00d796=\", \"88\u00d778=\"],\n  [\"43\u00d763=\", \"50\u00d768=\"],\n  [\"55\u00d723=\", \"73\u00d747=\"],\n  [\"28\u00d713=\", \"94\u00d749=\"],\n  [\"42\u00d761=\", \"37\u00d794=\"],\n  [\"81\u00d746=\", \"15\u00d715=\"],\n  [\"71\u00d744=\", \"33\u00d756=\"],\n  [\"82\u00d765=\", \"82\u00d783=\"],\n  [\"64\u00d783=\", \"48\u00d722=\"],\n  [\"20\u00d790=\", \"43\u00d789=\"],\n  [\"15\u00d735=\", \"68\u00d786=\"],\n  [\"76\u00d740=\", \"81\u00d752=\"],\n  [\"42\u00d731=\", \"47\u00d768=\"],\n  [\"27\u00d726=\", \"39\u00d771=\"],\n  [\"86\u00d737=\", \"69\u00d790=\"],\n  [\"77\u00d728=\", \"87\u00d727=\"],\n  [\"25\u00d736=\", \"54\u00d783=\"],\n  [\"36\u00d729=\", \"58\u00d794=\"],\n  [\"12\u00d793=\", \"41\u00d742=\"],\n  [\"78\u00d742=\", \"66\u00d740=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date heading and every two-digit-by-two-digit multiplication\n# prompt in the practice-sheet table. Each \"old\" value below occurs exactly\n# once in the document, so Find/Replace (wdReplaceAll) is safe and\n# unambiguous for every pair - it only ever touches its single match.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-26 Monday\", \"2024-02-27 Tuesday\"),\n    @(\"37\u00d788=\", \"66\u00d788=\"),\n    @(\"49\u00d723=\", \"40\u00d721=\"),\n    @(\"20\u00d739=\", \"99\u00d724=\"),\n    @(\"81\u00d774=\", \"37\u00d737=\"),\n    @(\"37\u00d773=\", \"12\u00d762=\"),\n    @(\"33\u00d796=\", \"88\u00d778=\"),\n    @(\"43\u00d763=\", \"50\u00d768=\"),\n    @(\"55\u00d723=\", \"73\u00d747=\"),\n    @(\"28\u00d713=\", \"94\u00d749=\"),\n    @(\"42\u00d761=\", \"37\u00d794=\"),\n    @(\"81\u00d746=\", \"15\u00d715=\"),\n    @(\"71\u00d744=\", \"33\u00d756=\"),\n    @(\"82\u00d765=\", \"82\u00d783=\"),\n    @(\"64\u00d783=\", \"48\u00d722=\"),\n    @(\"20\u00d790=\", \"43\u00d789=\"),\n    @(\"15\u00d735=\", \"68\u00d786=\"),\n    @(\"76\u00d740=\", \"81\u00d752=\"),\n    @(\"42\u00d731=\", \"47\u00d768=\"),\n    @(\"27\u00d726=\", \"39\u00d771=\"),\n    @(\"86\u00d737=\", \"69\u00d790=\"),\n    @(\"77\u00d728=\", \"87\u00d727=\"),\n    @(\"25\u00d736=\", \"54\u00d783=\"),\n    @(\"36\u00d729=\", \"58\u00d794=\"),\n    @(\"12\u00d793=\", \"41\u00d742=\"),\n    @(\"78\u00d742=\", \"66\u00d740=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #          MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #          Format, ReplaceWith, Replace)\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
